# Re-apply the new ordering of "goods" labels to column A (rows 2-56)
# while leaving the counts in column B untouched. This mirrors the
# upstream diff, which only reshuffled the shared-strings table (i.e.
# which label text is shown in which already-existing row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLabels = @(
    "хлеб",
    "вино",
    "скот",
    "холст",
    "кожа",
    "мед",
    "пиво",
    "сукно",
    "овчина",
    "лошадь",
    "воск",
    "масло",
    "сало",
    "железо",
    "полотно",
    "колеса",
    "Крымскую соль",
    "говядина",
    "сено",
    "парча",
    "позумент",
    "табак",
    "выбойка",
    "чулок",
    "сахар",
    "шелк",
    "лыко",
    "лес",
    "коса",
    "ладан",
    "китайка",
    "сани",
    "сапог",
    "замок",
    "обод",
    "веревка",
    "ром",
    "платок",
    "гвоздь",
    "рогожа",
    "горшок",
    "конь",
    "овца",
    "покроми",
    "бечева",
    "роза",
    "котел",
    "сосуд",
    "брусья",
    "хомут",
    "нитка",
    "дуга",
    "гумми",
    "скотский кожа",
    "сковорода"
)

for ($i = 0; $i -lt $newLabels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newLabels[$i]
}
